$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("91th") | Out-Null
$s = $r.Start
$e = $r.End
$leftChar = $d.Range($s - 1, $s).Text
$rightChar = $d.Range($e, $e + 1).Text
$big = $d.Range($s - 1, $e + 1)
$big.Text = $leftChar + "91st" + $rightChar
$r2 = $d.Content
$r2.Find.Execute("91st") | Out-Null
$sub = $d.Range($r2.Start, $r2.End)
$sub.Font.Engrave = 1
$sub.Font.Engrave = 0
Write-Output ("doc: " + $d.Content.Text.Substring($r2.Start-20, 60))
